$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 4
$ws.Range("E4").Value = 10
$ws.Range("F4").Value = 2
$ws.Range("H4").Value = 2

# Row 15
$ws.Range("E15").Value = 105

# Row 18
$ws.Range("E18").Value = 56
$ws.Range("F18").Value = 17
$ws.Range("H18").Value = 17

# Row 21
$ws.Range("E21").Value = 1

# Row 36
$ws.Range("E36").Value = 53

# Row 38
$ws.Range("E38").Value = 39

# Row 39
$ws.Range("E39").Value = 14

# Row 41
$ws.Range("E41").Value = 20

# Row 44
$ws.Range("E44").Value = 18
$ws.Range("F44").Value = 7
$ws.Range("H44").Value = 7

# Row 46
$ws.Range("E46").Value = 14
$ws.Range("F46").Value = 2
$ws.Range("H46").Value = 2

# Row 62
$ws.Range("E62").Value = 21

# Row 80
$ws.Range("E80").Value = 14
$ws.Range("F80").Value = 5
$ws.Range("H80").Value = 5
